$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.143.70"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "2.962.85"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.568"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.31%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.627"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.138"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0873"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("D14").Value = "3.461.39"
$ws.Range("E14").Value = "  +3.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "3.037.85"
$ws.Range("E16").Value = "  +4.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.985"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "52.184.16"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("E19").Value = "  +6.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "0.0₃0981"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.178"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.107"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.00%  "
$ws.Range("E34").Value = "  -2.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "52.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0441"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.81%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.119"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.39%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("D48").Value = "2.138.70"
$ws.Range("E48").Value = "  -1.35%  "
$ws.Range("E49").Value = "  +2.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.243"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.919"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.80%  "

